$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D17","E17","F17","G17","H17","I17","J17","K17","L17","M17","N17","O17","P17","Q17","R17","S17","T17","U17","V17","W17","X17","Y17","Z17","AA17")
$vals = @(5001,4813,4741,4727,4770,4884,5267,5685,6019,6120,6044,5977,5899,5847,5778,5768,5778,5887,6143,6160,5954,5697,5472,5120)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $vals[$i]
}

$cols = @("D18","E18","F18","G18","H18","I18","J18","K18","L18","M18","N18","O18","P18","Q18","R18","S18","T18","U18","V18","W18","X18","Y18","Z18","AA18")
$vals = @(5200,5012,4917,4885,4873,4931,5187,5645,6028,5990,5872,5789,5723,5697,5677,5731,5833,5987,6162,6077,5817,5604,5469,5238)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $vals[$i]
}

$cols = @("D19","E19","F19","G19","H19","I19","J19","K19","L19","M19","N19","O19","P19","Q19","R19","S19","T19","U19","V19","W19","X19","Y19","Z19","AA19")
$vals = @(5253,5097,5008,4951,4935,4968,5092,5341,5722,6017,6111,6114,6054,5986,5866,5810,5967,6187,6492,6495,6298,6029,5794,5445)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $vals[$i]
}

$cols = @("D20","E20","F20","G20","H20","I20","J20","K20","L20","M20","N20","O20","P20","Q20","R20","S20","T20","U20","V20","W20","X20","Y20","Z20","AA20")
$vals = @(5074,4882,4814,4791,4829,4984,5440,5902,6244,6374,6252,6214,6122,6078,6006,5995,6042,6180,6437,6477,6269,5990,5756,5376)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $vals[$i]
}

$cols = @("D21","E21","F21","G21","H21","I21","J21","K21","L21","M21","N21","O21","P21","Q21","R21","S21","T21","U21","V21","W21","X21","Y21","Z21","AA21")
$vals = @(4718,4542,4474,4460,4506,4628,5039,5488,5847,5955,5902,5855,5800,5764,5716,5709,5717,5810,6027,6041,5825,5554,5318,4948)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i]).Value = $vals[$i]
}

# New row 22
$ws.Range("A22").Value = 2026
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = 25
$cols22 = @("D22","E22","F22","G22","H22","I22","J22","K22","L22","M22","N22","O22","P22","Q22","R22","S22","T22","U22","V22","W22","X22","Y22","Z22","AA22")
$vals22 = @(4589,4413,4346,4332,4378,4500,4910,5358,5716,5825,5780,5741,5696,5666,5625,5619,5627,5706,5893,5905,5689,5419,5183,4814)
for ($i = 0; $i -lt $cols22.Length; $i++) {
    $ws.Range($cols22[$i]).Value = $vals22[$i]
}

Write-Output "done"